$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update Version / Date, replace the duplicated
#     "Contact" row with "Publisher" / "Jurisdiction" values, and drop the
#     now-redundant extra Contact row. ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (row 9) was blank -> "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was a duplicated "Contact" / "No display for ContactDetail" row;
# turn it into "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 duplicated the same "Contact" / "No display for ContactDetail"
# pair - remove it entirely, shifting everything below up by one row.
$meta.Rows.Item(11).Delete()

# --- Sheet "Elements": the Extension row's Short / Definition text now
#     describes the episode group code extension specifically. ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "Episode Group Code"
$elements.Range("L2").Value = "Group code for the episode of care"
